# Auto-generated update of D (Price) / E (Volume(1h)) columns for cryptos worksheet
# All target cells are plain text in the source data (inline strings), so we
# force the Text number format before assignment to avoid Excel auto-converting
# numeric-looking price strings (e.g. "613.22") into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.321.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.251.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.22"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.20"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.252.71"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.497"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.08"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.782.26"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.394.21"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.248.72"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "504.65"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.755"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.34"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.02"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.124"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +41.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.93"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.51"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +19.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0781"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +15.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "494.69"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0422"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.86"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.293"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.990.42"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.96"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.120"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.99%  "
